$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a text value into a cell while avoiding implicit numeric
# conversion (so strings like "2.83" stay text) and without leaving a
# leftover "text" number-format style on the cell afterwards.
# ---------------------------------------------------------------------------
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q1" worksheet by duplicating the most recent
#    quarter sheet ("2021-Q4"), which carries over all the sheet-level
#    formatting (page setup, header style, index-column style, ...).
#    The duplicate is placed right after it, i.e. right before "总计".
# ---------------------------------------------------------------------------
$prevQuarter = $wb.Worksheets.Item("2021-Q4")
$prevQuarter.Copy($null, $prevQuarter)

$newSheet = $wb.Worksheets.Item($prevQuarter.Index + 1)
$newSheet.Name = "2022-Q1"

# The source sheet only had one fund row (row 2); extend with two more rows
# (3 and 4), copying the same formatting (styled index cell in column A,
# unstyled data cells) down from row 2.
$newSheet.Range("A2:H2").Copy()
$newSheet.Range("A3:H4").PasteSpecial(-4122)

# Headers (B1:H1) - values only, formatting already came from the copy.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2
$newSheet.Range("A2").Value = 0
Set-TextValue $newSheet.Range("B2") "009613"
Set-TextValue $newSheet.Range("C2") "上银中证500指数增强A"
Set-TextValue $newSheet.Range("D2") "2.83"
Set-TextValue $newSheet.Range("E2") "90.41"
Set-TextValue $newSheet.Range("F2") "1.26"
Set-TextValue $newSheet.Range("G2") "0.0357"
$newSheet.Range("H2").Value = 5

# Row 3
$newSheet.Range("A3").Value = 1
Set-TextValue $newSheet.Range("B3") "009614"
Set-TextValue $newSheet.Range("C3") "上银中证500指数增强C"
Set-TextValue $newSheet.Range("D3") "1.70"
Set-TextValue $newSheet.Range("E3") "90.41"
Set-TextValue $newSheet.Range("F3") "1.26"
Set-TextValue $newSheet.Range("G3") "0.0214"
$newSheet.Range("H3").Value = 5

# Row 4
$newSheet.Range("A4").Value = 2
Set-TextValue $newSheet.Range("B4") "519165"
Set-TextValue $newSheet.Range("C4") "新华鑫利灵活配置混合"
Set-TextValue $newSheet.Range("D4") "0.05"
Set-TextValue $newSheet.Range("E4") "74.84"
Set-TextValue $newSheet.Range("F4") "3.42"
Set-TextValue $newSheet.Range("G4") "0.0017"
$newSheet.Range("H4").Value = 4

# ---------------------------------------------------------------------------
# 2. Update "总计" sheet: insert a new data row for 2022-Q1 at the top of the
#    data (row 2), pushing the existing rows down.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Range("A2:D2").Insert()
# Remove the auto-inherited style from the newly inserted (blank) row so the
# data cells end up with no explicit style, matching the other data rows.
$totalSheet.Range("B2:D2").Style = "Normal"
# Column A keeps the bold/bordered "index" style used by every other row;
# grab it from row 3 (the row right below, which still has its original
# style) since the freshly inserted A2 cell currently has no style at all.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.06

# Renumber the index column (A) for all the rows that got shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
